# Change in retailer name lookup
# Append two new retailer rows (City = Ahmedabad) to the lookup table on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 22: ALLIED AGENCY / Ahmedabad
$ws.Range("A22").Value = "ALLIED AGENCY"
$ws.Range("C22").Value = "Ahmedabad"

# New row 23: Aeroflon Engineers Pvt Ltd. / Ahmedabad
$ws.Range("A23").Value = "Aeroflon Engineers Pvt Ltd."
$ws.Range("C23").Value = "Ahmedabad"

# Match the row height used by the other data rows in the table.
$ws.Rows.Item(22).RowHeight = 13.8
$ws.Rows.Item(23).RowHeight = 13.8

# Leave the selection on the next empty row of the lookup column, as in the
# edited workbook.
[void]$ws.Range("C24").Select()
